$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temp solve of RWheel: set the Fitness column (C) for rows 2-12 to a
# constant value of 4058.
$ws.Range("C2:C12").Value = 4058
